# Apply Halicarnassus Profits leve-price/profit updates across all 8 class sheets.
# Values derived from scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 844.2308
$ws.Range("I92").Value = 184.375
$ws.Range("J92").Value = 1900
$ws.Range("K92").Value = 184.375
$ws.Range("L92").Value = 1900
$ws.Range("M92").Value = 1063.625
$ws.Range("N92").Value = -4396
$ws.Range("H96").Value = 605.4167
$ws.Range("I96").Value = 398.8889
$ws.Range("K96").Value = 1196.6667
$ws.Range("M96").Value = 176.3333
$ws.Range("H121").Value = 878.25
$ws.Range("J121").Value = 878.25
$ws.Range("L121").Value = 2634.75
$ws.Range("N121").Value = -6128.75
$ws.Range("H137").Value = 2697.2173
$ws.Range("I137").Value = 2056.0715
$ws.Range("J137").Value = 3694.5557
$ws.Range("K137").Value = 6168.2145
$ws.Range("L137").Value = 11083.6671
$ws.Range("M137").Value = -3618.2145
$ws.Range("N137").Value = -16183.6671

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2369.8
$ws.Range("J88").Value = 2462.25
$ws.Range("L88").Value = 2462.25
$ws.Range("N88").Value = -3274.25
$ws.Range("H91").Value = 2369.8
$ws.Range("J91").Value = 2462.25
$ws.Range("L91").Value = 2462.25
$ws.Range("N91").Value = -5270.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7566.8667
$ws.Range("J20").Value = 4926.5
$ws.Range("L20").Value = 4926.5
$ws.Range("N20").Value = -5420.5
$ws.Range("H63").Value = 75000
$ws.Range("J63").Value = 75000
$ws.Range("L63").Value = 75000
$ws.Range("N63").Value = -76372
$ws.Range("H66").Value = 75000
$ws.Range("J66").Value = 75000
$ws.Range("L66").Value = 225000
$ws.Range("N66").Value = -231864
$ws.Range("H86").Value = 7655.25
$ws.Range("I86").Value = 622
$ws.Range("K86").Value = 622
$ws.Range("M86").Value = 501
$ws.Range("H89").Value = 7655.25
$ws.Range("I89").Value = 622
$ws.Range("K89").Value = 3110
$ws.Range("M89").Value = 2506
$ws.Range("H107").Value = 4366.1
$ws.Range("I107").Value = 1360.1666
$ws.Range("K107").Value = 1360.1666
$ws.Range("M107").Value = 559.8334
$ws.Range("H127").Value = 57443.5
$ws.Range("I127").Value = 54888
$ws.Range("J127").Value = 59999
$ws.Range("K127").Value = 54888
$ws.Range("L127").Value = 59999
$ws.Range("M127").Value = -49928
$ws.Range("N127").Value = -69919

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6686.8335
$ws.Range("I31").Value = 1578.6666
$ws.Range("J31").Value = 9751.733
$ws.Range("K31").Value = 1578.6666
$ws.Range("L31").Value = 9751.733
$ws.Range("M31").Value = -1283.6666
$ws.Range("N31").Value = -10341.733
$ws.Range("H34").Value = 6686.8335
$ws.Range("I34").Value = 1578.6666
$ws.Range("J34").Value = 9751.733
$ws.Range("K34").Value = 1578.6666
$ws.Range("L34").Value = 9751.733
$ws.Range("M34").Value = -1376.6666
$ws.Range("N34").Value = -10155.733
$ws.Range("H58").Value = 2349.8823
$ws.Range("I58").Value = 1421.1
$ws.Range("K58").Value = 1421.1
$ws.Range("M58").Value = -1218.1
$ws.Range("H122").Value = 1464.15
$ws.Range("I122").Value = 1313.1765
$ws.Range("K122").Value = 3939.5295
$ws.Range("M122").Value = -1489.5295
$ws.Range("H136").Value = 2349.8823
$ws.Range("I136").Value = 1421.1
$ws.Range("K136").Value = 4263.299999999999
$ws.Range("M136").Value = -1713.299999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1106.1538
$ws.Range("I5").Value = 1041.5714
$ws.Range("K5").Value = 3124.7142
$ws.Range("M5").Value = -3012.7142
$ws.Range("H52").Value = 1000
$ws.Range("I52").Value = 1000
$ws.Range("K52").Value = 3000
$ws.Range("M52").Value = -2734
$ws.Range("H55").Value = 4581.25
$ws.Range("I55").Value = 1162.5
$ws.Range("K55").Value = 3487.5
$ws.Range("M55").Value = -3310.5
$ws.Range("H132").Value = 2308.6365
$ws.Range("I132").Value = 2059.8
$ws.Range("K132").Value = 18538.2
$ws.Range("M132").Value = -16008.2
$ws.Range("H135").Value = 1106.1538
$ws.Range("I135").Value = 1041.5714
$ws.Range("K135").Value = 9374.142600000001
$ws.Range("M135").Value = -6839.142600000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("N73").ClearContents()
$ws.Range("H80").Value = 2062.125
$ws.Range("J80").Value = 2249.5
$ws.Range("L80").Value = 2249.5
$ws.Range("N80").Value = -4245.5
$ws.Range("H83").Value = 2062.125
$ws.Range("J83").Value = 2249.5
$ws.Range("L83").Value = 11247.5
$ws.Range("N83").Value = -21231.5
$ws.Range("H132").Value = 15976.75
$ws.Range("I132").Value = 17561.445
$ws.Range("J132").Value = 11222.667
$ws.Range("K132").Value = 52684.335
$ws.Range("L132").Value = 33668.001
$ws.Range("M132").Value = -50154.335
$ws.Range("N132").Value = -38728.001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1078.1666
$ws.Range("I22").Value = 1183.8
$ws.Range("K22").Value = 1183.8
$ws.Range("M22").Value = -888.8
$ws.Range("H27").Value = 1078.1666
$ws.Range("I27").Value = 1183.8
$ws.Range("K27").Value = 1183.8
$ws.Range("M27").Value = -1076.8
$ws.Range("H46").Value = 3976.8096
$ws.Range("I46").Value = 2049
$ws.Range("J46").Value = 5422.6665
$ws.Range("K46").Value = 2049
$ws.Range("L46").Value = 5422.6665
$ws.Range("M46").Value = -1861
$ws.Range("N46").Value = -5798.6665
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("H82").Value = 3302.8
$ws.Range("I82").Value = 2299.3333
$ws.Range("K82").Value = 2299.3333
$ws.Range("M82").Value = -1938.3333
$ws.Range("H85").Value = 3302.8
$ws.Range("I85").Value = 2299.3333
$ws.Range("K85").Value = 2299.3333
$ws.Range("M85").Value = -1051.3333
$ws.Range("H122").Value = 4039.6
$ws.Range("I122").Value = 3799.5
$ws.Range("K122").Value = 11398.5
$ws.Range("M122").Value = -8948.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H112").Value = 31258
$ws.Range("J112").Value = 31258
$ws.Range("L112").Value = 31258
$ws.Range("N112").Value = -34212
$ws.Range("H136").Value = 2681.457
$ws.Range("I136").Value = 1666.4445
$ws.Range("K136").Value = 4999.333500000001
$ws.Range("M136").Value = -2449.333500000001
$ws.Range("H140").Value = 42900
$ws.Range("J140").Value = 42900
$ws.Range("L140").Value = 42900
$ws.Range("N140").Value = -53260
